$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-01 Tuesday" "2025-07-02 Wednesday"

Replace-Text "661÷7=94, 3" "535÷6=89, 1"
Replace-Text "354÷2=177, 0" "387÷4=96, 3"
Replace-Text "171÷7=24, 3" "406÷9=45, 1"
Replace-Text "424÷4=106, 0" "370÷7=52, 6"
Replace-Text "683÷7=97, 4" "895÷9=99, 4"

Replace-Text "171÷8=21, 3" "850÷4=212, 2"
Replace-Text "650÷2=325, 0" "376÷4=94, 0"
Replace-Text "184÷3=61, 1" "546÷2=273, 0"
Replace-Text "236÷6=39, 2" "568÷8=71, 0"
Replace-Text "149÷3=49, 2" "385÷5=77, 0"

Replace-Text "987÷4=246, 3" "537÷5=107, 2"
Replace-Text "800÷7=114, 2" "791÷4=197, 3"
Replace-Text "997÷9=110, 7" "521÷4=130, 1"
Replace-Text "858÷2=429, 0" "161÷5=32, 1"
Replace-Text "584÷7=83, 3" "698÷2=349, 0"

Replace-Text "869÷8=108, 5" "103÷8=12, 7"
Replace-Text "223÷6=37, 1" "771÷5=154, 1"
Replace-Text "348÷8=43, 4" "441÷8=55, 1"
Replace-Text "442÷4=110, 2" "180÷7=25, 5"
Replace-Text "404÷9=44, 8" "480÷2=240, 0"

Replace-Text "764÷2=382, 0" "438÷5=87, 3"
Replace-Text "659÷7=94, 1" "272÷5=54, 2"
Replace-Text "987÷8=123, 3" "973÷5=194, 3"
Replace-Text "736÷9=81, 7" "445÷8=55, 5"
Replace-Text "202÷3=67, 1" "214÷9=23, 7"
